$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after the first sheet (so it becomes sheet index 2)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "SendUserMessage-Event"

# Header row
$newSheet.Range("A1").Value = "TestCaseName"
$newSheet.Range("B1").Value = "Type"
$newSheet.Range("C1").Value = "Resource"
$newSheet.Range("D1").Value = "TestCaseNameDesc"
$newSheet.Range("E1").Value = "StepInfo"
$newSheet.Range("F1").Value = "Event"
$newSheet.Range("G1").Value = "Identifier"
$newSheet.Range("H1").Value = "RequestContent"
$newSheet.Range("I1").Value = "MessageType"
$newSheet.Range("J1").Value = "Csvson"
$newSheet.Range("K1").Value = "Tags"

# Row 2
$newSheet.Range("A2").Value = "SEND_USER"
$newSheet.Range("B2").Value = "KAFKA"
$newSheet.Range("C2").Value = "avro"
$newSheet.Range("D2").Value = "Create User "
$newSheet.Range("E2").Value = "As a user needs to;create user;details;user;user;contains user information in api;"
$newSheet.Range("F2").Value = "UserCreated"

# Row 3
$newSheet.Range("A3").Value = "VERIFY_USER_CREATED_EVENT_1"
$newSheet.Range("B3").Value = "KAFKA"
$newSheet.Range("C3").Value = "avro"
$newSheet.Range("D3").Value = "contains user information"
$newSheet.Range("E3").Value = "Validate created user event"
$newSheet.Range("F3").Value = "UserCreated"

# Remaining row 2 cells
$newSheet.Range("K2").Value = " @simple-kafka-send"
$newSheet.Range("I2").Value = "AvroType"

# Remaining row 3 cells
$newSheet.Range("G3").Value = "Rockey"
$newSheet.Range("J3").Value = "name,age`nRockey,i~44"
$newSheet.Range("I3").Value = "AvroType"
$newSheet.Range("K3").Value = "  @validate_kafka_message @IncludesByPath"

# Last new string
$newSheet.Range("H2").Value = "{ ""name"" : ""Rocky"",""age"":`n44 }"
